$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order for the data rows below: B, D, E, F, G, I, K, N
$cols = @(2, 4, 5, 6, 7, 9, 11, 14)

$data = @(
    @(0.04899958594673137, 0.04890864427885688, 0.4099742985248014, 1.647936557103208, 0.002458512755508271, 0.5535243652359512, 1.779141658815718, 1.598384161986061),
    @(0.04336700276643768, 0.0492964254672188, 0.3573310487299324, 1.594934397641211, 0.002464529494027036, 0.5556605943118704, 1.594978902241678, 1.612384793975117),
    @(0.03991087294078, 0.04956484607794565, 0.32515709915522, 1.563605916047322, 0.002468409692499588, 0.5572233815539107, 1.482796537208685, 1.621578194840325),
    @(0.03850321301737836, 0.04968177236794702, 0.3120801918104661, 1.551141417440363, 0.002470037839598734, 0.5579231827406801, 1.437299830781797, 1.625474100240524),
    @(0.03826952172735787, 0.04970164108434716, 0.3099107541411001, 1.549089844600928, 0.002470311032634097, 0.5580431806574282, 1.429758151441149, 1.626130025885736),
    @(0.03989188552924361, 0.04956639255313888, 0.3249806050883137, 1.563436596702445, 0.002468431459995206, 0.5572325646809766, 1.482182076276274, 1.621630131750344),
    @(0.04705708974218226, 0.04903600831773502, 0.3917899312540243, 1.62940713192819, 0.002460548858358967, 0.5542086944914075, 1.715452808917235, 1.603087325728907),
    @(0.06111983681967104, 0.04824026631132838, 0.5241424050174288, 1.768580737753837, 0.002446557421123155, 0.5502807600826323, 2.18030609276019, 1.571487949218174),
    @(0.07145133554138283, 0.04780988175598821, 0.6224248602630951, 1.877052154206154, 0.002437159600223961, 0.5486285883411952, 2.526825505937609, 1.551211933603668),
    @(0.076149782467823, 0.0476487622511641, 0.6674130776613225, 1.927803640628525, 0.002433073108738219, 0.5481478054799851, 2.685656841289244, 1.54263429102231),
    @(0.07792861126658579, 0.04759283544316162, 0.6844931037474993, 1.947228551039558, 0.002431552581094576, 0.5480049266515721, 2.745982573788353, 1.539479759248351),
    @(0.07754552784359703, 0.04760465246238965, 0.6808125971816423, 1.943035802766843, 0.002431878858976522, 0.5480339515702113, 2.732982232186657, 1.540154967726636),
    @(0.07629613618523479, 0.04764405860796472, 0.668817358277181, 1.929397579123275, 0.002432947475163047, 0.5481352640592547, 2.690616229544389, 1.542372884924134),
    @(0.07553079439095711, 0.04766886135131188, 0.6614757689624327, 1.921070777655018, 0.002433605536249005, 0.5482024309635705, 2.664689464157163, 1.543743639490202),
    @(0.0711442382462053, 0.04782111740414052, 0.6194907156947664, 1.873764088466373, 0.00243743044150748, 0.5486654784473757, 2.51647024644052, 1.551785565077267),
    @(0.06845274782799038, 0.04792347100906369, 0.5938084190187709, 1.845106269159459, 0.002439825076206287, 0.5490190725099424, 2.425854649977055, 1.556885007822345),
    @(0.06690455164824982, 0.047985598704404, 0.5790626707886446, 1.828755263178806, 0.002441220172074949, 0.549247924546421, 2.373847250533345, 1.559878833736825),
    @(0.06638034323255226, 0.04800719036507317, 0.5740743687862135, 1.823241679059407, 0.002441695584437498, 0.5493297777733019, 2.356257499945059, 1.560902905272528),
    @(0.06873927544239677, 0.04791223752341978, 0.5965396218323633, 1.848143230212685, 0.002439568325918149, 0.5489787936025223, 2.43548915907752, 1.556335869568812),
    @(0.07666312421523003, 0.04763234522815907, 0.672339426516956, 1.933397820969248, 0.002432632867204777, 0.5481044407498885, 2.703055212458423, 1.541718881684645),
    @(0.08183957679544562, 0.04747911114676384, 0.7221370817743917, 1.99032176387135, 0.002428257072606777, 0.5477615057918825, 2.87897643635506, 1.532711966407646),
    @(0.07907706726611252, 0.04755814425668348, 0.6955342292189073, 1.959828751652026, 0.002430578218926072, 0.5479235471582413, 2.784985248938654, 1.537468907733412),
    @(0.06860973885676458, 0.04791730596341637, 0.5953047850704536, 1.846769831838913, 0.002439684345419355, 0.5489969240822603, 2.431133121753021, 1.556583941637655),
    @(0.05731501751981227, 0.04842883869280357, 0.4881709990556828, 1.729856005680574, 0.002450186758937369, 0.5511276956902904, 2.053709198861611, 1.579522791424509)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($startRow + $i, $cols[$j]).Value = $rowVals[$j]
    }
}
